$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap the "Periodo Mora" values between the two data rows
$ws.Range("E16").Value = "1802"
$ws.Range("E17").Value = "1803"

# Update "Valor Mora" amounts for both rows
$ws.Range("G16").Value = 737717
$ws.Range("G17").Value = 737717
